$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Test Plan Author:Avital Z" -> "Test Plan Author: Avital Z"
#    Also merges the "Author" bold run with "Test Plan " and strips the
#    spell/grammar proof marks that used to bracket "Author:Avital".
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Test Plan Author:Avital Z", $true, $false, $false, $false, $false, $true, 1, $false, "Test Plan Author: Avital Z", 2)

$rng = $d.Content
$null = $rng.Find.Execute("Test Plan Author: Avital Z", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
# "Test Plan Author" stays bold; everything from the colon onward is not bold.
$tail = $d.Range($start + 16, $end)
$tail.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) API Testing bullet: add "login authorization, authentication, " before
#    "order placement, order status updates, and MongoDB interaction."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Validate the REST API for order placement, order status updates, and MongoDB interaction.", $true, $false, $false, $false, $false, $true, 1, $false, "Validate the REST API for login authorization, authentication, order placement, order status updates, and MongoDB interaction.", 2)

# ---------------------------------------------------------------------------
# 3) "Verify email notification delivery." -> "...delivery to the customer
#    for each status change." and pick up the missing en-US language tag.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Verify email notification delivery.", $true, $false, $false, $false, $false, $true, 1, $false, "Verify email notification delivery to the customer for each status change.", 2)

$rng2 = $d.Content
$null = $rng2.Find.Execute("Verify email notification delivery to the customer for each status change.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.LanguageID = "en-US"

# ---------------------------------------------------------------------------
# 4) CRUD Operations bullet: add "login authorization, " before the
#    "order creation, fetching, updating, and deletion." text.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("CRUD Operations on orders: Test order creation, fetching, updating, and deletion.", $true, $false, $false, $false, $false, $true, 1, $false, "CRUD Operations on orders: Test login authorization, order creation, fetching, updating, and deletion.", 2)

# ---------------------------------------------------------------------------
# 5) "All bugs have been fixed..." -> "Critical high priority bugs have
#    been fixed..."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("All bugs have been fixed and verified in the test environment.", $true, $false, $false, $false, $false, $true, 1, $false, "Critical high priority bugs have been fixed and verified in the test environment.", 2)

# ---------------------------------------------------------------------------
# 6) Rendering bookkeeping: the stale <w:lastRenderedPageBreak/> markers that
#    used to sit in front of "Go to the URL ..." and "Login Flow: ..." no
#    longer belong there once the document reflows (the break now falls at
#    "Version:" / "Verifying UI Updates ..." instead); re-touching that text
#    drops the stale marker.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Go to the URL http://localhost:3000/user", $true, $false, $false, $false, $false, $true, 1, $false, "Go to the URL http://localhost:3000/user", 2)
$null = $d.Content.Find.Execute("Login Flow: Test login for multiple user roles (customer, admin). ", $true, $false, $false, $false, $false, $true, 1, $false, "Login Flow: Test login for multiple user roles (customer, admin). ", 2)
